$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) from the latest scrape.
# The "Price" column (D) stores numeric-looking values as plain text in the source
# (e.g. "66.814.28" with dotted separators, or "1.00" with a significant trailing
# zero). A leading apostrophe forces Excel to keep each one as literal text - matching
# the workbook's original inlineStr cells - instead of auto-coercing it into a Double
# and silently losing the separators / trailing zeros.

# Row 2
$ws.Range('D2').Value = "'66.831.42"
$ws.Range('E2').Value = '  +1.34%  '

# Row 3
$ws.Range('D3').Value = "'3.794.99"
$ws.Range('E3').Value = '  -0.50%  '

# Row 4
$ws.Range('E4').Value = '  -0.24%  '

# Row 5
$ws.Range('D5').Value = "'443.03"
$ws.Range('E5').Value = '  +5.63%  '

# Row 6
$ws.Range('D6').Value = "'144.50"
$ws.Range('E6').Value = '  +13.88%  '

# Row 7
$ws.Range('D7').Value = "'0.620"
$ws.Range('E7').Value = '  +3.41%  '

# Row 9
$ws.Range('D9').Value = "'0.732"
$ws.Range('E9').Value = '  +3.00%  '

# Row 10
$ws.Range('E10').Value = '  -6.96%  '

# Row 11
$ws.Range('D11').Value = "'0.0000309"
$ws.Range('E11').Value = '  -9.17%  '

# Row 12
$ws.Range('D12').Value = "'43.44"
$ws.Range('E12').Value = '  +9.15%  '

# Row 13
$ws.Range('D13').Value = "'10.32"
$ws.Range('E13').Value = '  +5.21%  '

# Row 14
$ws.Range('D14').Value = "'4.395.64"
$ws.Range('E14').Value = '  -0.92%  '

# Row 15
$ws.Range('D15').Value = "'14.71"
$ws.Range('E15').Value = '  -6.91%  '

# Row 16
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = "'0.137"
$ws.Range('E16').Value = '  -0.16%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = "'3.764.78"
$ws.Range('E17').Value = '  -1.20%  '

# Row 18
$ws.Range('D18').Value = "'19.82"
$ws.Range('E18').Value = '  +2.61%  '

# Row 19
$ws.Range('D19').Value = "'1.14"
$ws.Range('E19').Value = '  +7.88%  '

# Row 20
$ws.Range('D20').Value = "'66.835.81"
$ws.Range('E20').Value = '  +0.92%  '

# Row 21
$ws.Range('D21').Value = "'414.37"
$ws.Range('E21').Value = '  +3.78%  '

# Row 22
$ws.Range('D22').Value = "'14.52"
$ws.Range('E22').Value = '  +3.05%  '

# Row 23
$ws.Range('E23').Value = '  +10.71%  '

# Row 24
$ws.Range('D24').Value = "'85.44"
$ws.Range('E24').Value = '  +2.68%  '

# Row 25
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').Value = "'36.94"
$ws.Range('E25').Value = '  +0.76%  '

# Row 26
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = "'3.40"
$ws.Range('E26').Value = '  +7.81%  '

# Row 27
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').Value = "'5.52"
$ws.Range('E27').Value = '  -4.12%  '

# Row 28
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = "'9.55"
$ws.Range('E28').Value = '  +29.09%  '

# Row 29
$ws.Range('D29').Value = "'9.72"
$ws.Range('E29').Value = '  +4.59%  '

# Row 30
$ws.Range('D30').Value = "'732.86"
$ws.Range('E30').Value = '  +5.17%  '

# Row 31
$ws.Range('D31').Value = "'13.84"
$ws.Range('E31').Value = '  +13.86%  '

# Row 32
$ws.Range('E32').Value = '  +11.97%  '

# Row 33
$ws.Range('D33').Value = "'2.74"
$ws.Range('E33').Value = '  -0.13%  '

# Row 34
$ws.Range('D34').Value = "'43.69"
$ws.Range('E34').Value = '  +16.88%  '

# Row 35
$ws.Range('D35').Value = "'0.159"
$ws.Range('E35').Value = '  +7.15%  '

# Row 36
$ws.Range('D36').Value = "'56.52"
$ws.Range('E36').Value = '  +3.57%  '

# Row 37
$ws.Range('E37').Value = '  +0.13%  '

# Row 38
$ws.Range('D38').Value = "'5.49"
$ws.Range('E38').Value = '  +25.40%  '

# Row 39
$ws.Range('D39').Value = "'0.0476"
$ws.Range('E39').Value = '  +6.58%  '

# Row 40
$ws.Range('D40').Value = "'2.87"
$ws.Range('E40').Value = '  -0.29%  '

# Row 41
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = "'0.340"
$ws.Range('E41').Value = '  +19.31%  '

# Row 42
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = "'2.59"
$ws.Range('E42').Value = '  +30.69%  '

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = "'0.141"
$ws.Range('E43').Value = '  +5.52%  '

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = "'1.00"
$ws.Range('E44').Value = '  -0.22%  '

# Row 45
$ws.Range('D45').Value = "'0.0₃0671"
$ws.Range('E45').Value = '  -9.80%  '

# Row 46
$ws.Range('D46').Value = "'3.27"
$ws.Range('E46').Value = '  +7.02%  '

# Row 47
$ws.Range('D47').Value = "'3.33"
$ws.Range('E47').Value = '  +1.64%  '

# Row 48
$ws.Range('D48').Value = "'144.91"
$ws.Range('E48').Value = '  +1.52%  '

# Row 49
$ws.Range('D49').Value = "'2.09"
$ws.Range('E49').Value = '  +2.70%  '

# Row 50
$ws.Range('E50').Value = '  +5.01%  '

# Row 51
$ws.Range('D51').Value = "'2.83"
$ws.Range('E51').Value = '  +5.02%  '
